# Applies the PENNSYLVANIA_2016.xlsx cleanup:
#  1. Renames the header row to short machine-friendly column names.
#  2. Title-cases Spanish connector words (de/del/el/la/los/las/y) inside the
#     "Estado de Origen" (A) and "Municipio Origen" (B) text columns.
#  3. Fixes a floating point value that was stored with excess precision
#     (5/5389) so it matches Excel's normal 15-significant-digit rounding.
#  4. Removes the trailing footnote/source rows (937-941), which shrinks the
#     used range down to A1:D935.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Header row -----------------------------------------------------
$ws.Range("A1").Value = "mx_state"
$ws.Range("B1").Value = "mx_municipality"
$ws.Range("C1").Value = "n_matriculas"
$ws.Range("D1").Value = "pct_matriculas"

# --- 2. Title-case connector words in columns A and B ------------------
$lastRow = 935
for ($r = 2; $r -le $lastRow; $r++) {
    foreach ($col in 1, 2) {
        $cell = $ws.Cells.Item($r, $col)
        $v = $cell.Value2
        if (($v -ne $null) -and ($v -is [string])) {
            $nv = $v -replace '\bde\b','De' -replace '\bdel\b','Del' -replace '\bel\b','El' -replace '\bla\b','La' -replace '\blos\b','Los' -replace '\blas\b','Las' -replace '\by\b','Y'
            $cell.Value = $nv
        }
    }
}

# --- 3. Fix the over-precise percentage literal (5/5389) ---------------
for ($r = 2; $r -le $lastRow; $r++) {
    $dcell = $ws.Cells.Item($r, 4)
    $dv = $dcell.Value2
    if (($dv -ne $null) -and (-not ($dv -is [string]))) {
        if ($dv -eq 0.0009278159213212099) {
            $dcell.Value = 0.00092781592132121
        }
    }
}

# --- 4. Drop the trailing footnote rows (937-941) -----------------------
$ws.Range("A937:A941").EntireRow.Delete()
